$d = $word.ActiveDocument

# --- 1) Header contact block: turn plain-text contact strings into hyperlinks ---

# email
$r = $d.Content
$r.Find.Execute("kaleb.coberly@gmail.com", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$d.Hyperlinks.Add($r, "mailto:kaleb.coberly@gmail.com") | Out-Null

# linkedin - split across "linkedin.com/in/" and "kaleb-coberly" (the latter wrapped in proofErr)
$r = $d.Content
$r.Find.Execute("linkedin.com/in/", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$d.Hyperlinks.Add($r, "https://linkedin.com/in/kaleb-coberly") | Out-Null

$r = $d.Content
$r.Find.Execute("kaleb-coberly", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$d.Hyperlinks.Add($r, "https://linkedin.com/in/kaleb-coberly") | Out-Null

# github - split across "github.com/" and "KalebCoberly" (the latter wrapped in proofErr)
$r = $d.Content
$r.Find.Execute("github.com/", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$d.Hyperlinks.Add($r, "https://github.com/KalebCoberly") | Out-Null

$r = $d.Content
$r.Find.Execute("KalebCoberly", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$d.Hyperlinks.Add($r, "https://github.com/KalebCoberly") | Out-Null

# personal website
$r = $d.Content
$r.Find.Execute("kalebcoberly.com", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$d.Hyperlinks.Add($r, "https://kalebcoberly.com") | Out-Null

# --- 2) Shorten the displayed "cricketsandcomb.org/#projects" link text (drop the "https://" prefix) ---

$r = $d.Content
$r.Find.Execute("https://cricketsandcomb.org/#projects", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$prefix = $d.Range($r.Start, $r.Start + 8)
$prefix.Delete()

# --- 3) Education: "Data Management" -> "Database Management" ---

$r = $d.Content
$r.Find.Execute("B.S., Data", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$insertAt = $d.Range($r.End, $r.End)
$insertAt.InsertAfter("base")
